$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "26.773.91"
$ws.Cells.Item(2, 5).Value = "  -0.81%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.797.15"
$ws.Cells.Item(3, 5).Value = "  -1.06%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.9998"
$ws.Cells.Item(4, 5).Value = "  -0.11%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "308.63"
$ws.Cells.Item(5, 5).Value = "  -0.65%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.9993"
$ws.Cells.Item(6, 5).Value = "  -0.13%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4394"
$ws.Cells.Item(7, 5).Value = "  +4.17%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3678"
$ws.Cells.Item(8, 5).Value = "  +0.33%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.07347"
$ws.Cells.Item(9, 5).Value = "  +1.96%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.8559"
$ws.Cells.Item(10, 5).Value = "  +1.99%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "20.62"
$ws.Cells.Item(11, 5).Value = "  -0.54%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.953.66"
$ws.Cells.Item(12, 5).Value = "  +7.66%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "6.620"
$ws.Cells.Item(13, 5).Value = "  -0.13%  "

$ws.Cells.Item(14, 2).Value = "Litecoin"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "91.97"
$ws.Cells.Item(14, 5).Value = "  +3.39%  "

$ws.Cells.Item(15, 2).Value = "TRON"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.07069"
$ws.Cells.Item(15, 5).Value = "  +0.09%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "5.265"
$ws.Cells.Item(16, 5).Value = "  -0.22%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "1.000"
$ws.Cells.Item(17, 5).Value = "  -0.20%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.000008635"
$ws.Cells.Item(18, 5).Value = "  -1.56%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.9999"
$ws.Cells.Item(19, 5).Value = "  -0.10%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "14.76"
$ws.Cells.Item(20, 5).Value = "  -1.07%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "26.812.63"
$ws.Cells.Item(21, 5).Value = "  -0.86%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.145"
$ws.Cells.Item(22, 5).Value = "  +0.63%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "10.81"
$ws.Cells.Item(23, 5).Value = "  +0.06%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "1.973"
$ws.Cells.Item(24, 5).Value = "  -0.05%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "151.42"
$ws.Cells.Item(25, 5).Value = "  -0.09%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.200"
$ws.Cells.Item(26, 5).Value = "  -0.81%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "18.35"
$ws.Cells.Item(27, 5).Value = "  +0.62%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "5.170"
$ws.Cells.Item(28, 5).Value = "  -0.63%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "117.14"
$ws.Cells.Item(29, 5).Value = "  +1.11%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.08789"
$ws.Cells.Item(30, 5).Value = "  +0.43%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.7374"
$ws.Cells.Item(31, 5).Value = "  +0.09%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.153"
$ws.Cells.Item(32, 5).Value = "  -1.54%  "

$ws.Cells.Item(33, 2).Value = "Filecoin"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "4.428"
$ws.Cells.Item(33, 5).Value = "  +0.61%  "

$ws.Cells.Item(34, 2).Value = "HuobiToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "2.896"
$ws.Cells.Item(34, 5).Value = "  -2.17%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.9987"
$ws.Cells.Item(35, 5).Value = "  -0.18%  "

$ws.Cells.Item(36, 5).Value = "  -0.45%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.01952"
$ws.Cells.Item(37, 5).Value = "  +0.09%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.05176"
$ws.Cells.Item(38, 5).Value = "  -0.87%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.5231"
$ws.Cells.Item(39, 5).Value = "  +4.22%  "

$ws.Cells.Item(40, 2).Value = "MXToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.821"
$ws.Cells.Item(40, 5).Value = "  -1.45%  "

$ws.Cells.Item(41, 2).Value = "FraxShare"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "6.994"
$ws.Cells.Item(41, 5).Value = "  -3.75%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.1681"
$ws.Cells.Item(42, 5).Value = "  -0.11%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.5039"
$ws.Cells.Item(43, 5).Value = "  +7.12%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "8.447"
$ws.Cells.Item(44, 5).Value = "  -1.68%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.976"
$ws.Cells.Item(45, 5).Value = "  +5.75%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "10.39"
$ws.Cells.Item(46, 5).Value = "  -0.73%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "104.97"
$ws.Cells.Item(47, 5).Value = "  -1.00%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.9981"
$ws.Cells.Item(48, 5).Value = "  -0.21%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.659"
$ws.Cells.Item(49, 5).Value = "  +1.18%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.06282"
$ws.Cells.Item(50, 5).Value = "  -0.92%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.9131"
$ws.Cells.Item(51, 5).Value = "  +1.58%  "
